# Graded Discussion Topic 10 - apply the authored edits.
$d = $word.ActiveDocument

# 1) Remove the leftover "_GoBack" bookmark (bookmarkStart/bookmarkEnd pair)
#    that trailed the title run.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Expand the "in a big scale ..." sentence into the fuller Boeing /
#    computer-vision paragraph.
$old1 = "in a big scale I would like to move my career towards this amazing and growing field. "
$new1 = "in a big scale I would like to start exploring the current capabilities of the industry and what the current research is providing on this field. At Boeing there is a lot of interest on computer vision and some of the possible applications at the early stages of the manufacturing process for new airplanes and of course all the other applications for the space sector."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 3) Reword the "Apart from working on some personal projects..." sentence.
$old2 = "Apart from working on some personal projects to apply the knowledge that I got from this certificate there are some classes/certifications I am planning to go to. The following are just some of the ones that I have scheduled for the upcoming years:"
$new2 = "Apart from working on some personal projects to apply some of what I learned throughout this class/certification, there are some classes/certifications I am planning to attend to. The following are just some of the ones that I have scheduled for the upcoming couple of years:"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
